$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Pais" (country) column with validated country values for each record
$ws.Range("D1").Value = "Pais"
$ws.Range("D2").Value = "Panama"
$ws.Range("D3").Value = "Colombia"
$ws.Range("D4").Value = "Colombia"
$ws.Range("D5").Value = "Panama"

# Narrow column C a bit to make room for the new column / improve the view
$ws.Columns.Item(3).ColumnWidth = 28.15

# Match the author's final selection/view state
[void]$ws.Range("D6").Select()
